# Added test case number 35
$wb = $excel.ActiveWorkbook

$wsBuySell = $wb.Worksheets.Item("NitroXBuySell")
$wsBots    = $wb.Worksheets.Item("NitroXBots")

# --- NitroXBuySell (sheet2): fill in the "OpenOrderNumber"/"SkipAtStepNum" cells on row 17 ---
$wsBuySell.Range("Q17").Value = 1
$wsBuySell.Range("R17").Value = 1

# --- Row 18: add SkipAtStepNum value ---
$wsBuySell.Range("R18").Value = 1

# --- Row 19: new test case QA_TestCase_Auto_NitroX_035, copied from row 18 with new TestCaseID ---
$wsBuySell.Range("A19").Value = "QA_TestCase_Auto_NitroX_035"
$wsBuySell.Range("B19").Value = "Futures"
$wsBuySell.Range("C19").Value = "QUANT_BINANCEDM_25"
$wsBuySell.Range("F19").Value = "ETH/USDT Perpetual USDT"
$wsBuySell.Range("H19").Value = "One Way"
$wsBuySell.Range("I19").Value = 25
$wsBuySell.Range("J19").Value = "CROSSED"
$wsBuySell.Range("K19").Value = "BUY"
$wsBuySell.Range("N19").Value = 1
$wsBuySell.Range("Q19").Value = 1
$wsBuySell.Range("R19").Value = 1
$wsBuySell.Range("S19").Value = $false

# --- Selection / active sheet bookkeeping ---
$wsBots.Activate()
$wsBots.Range("K15").Select()

$wsBuySell.Activate()
$wsBuySell.Range("A19").Select()

$wsBots.Activate()
